# Insert a new row at row 157 (shifts existing rows 157-229 down to 158-230)
# and populate it with a new Membrillo price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("157:157").Insert()

$ws.Cells.Item(157, 1).Value  = 10
$ws.Cells.Item(157, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(157, 3).Value  = "La Araucanía"
$ws.Cells.Item(157, 4).Value  = 44813
$ws.Cells.Item(157, 5).Value  = 9
$ws.Cells.Item(157, 6).Value  = "Fruta"
$ws.Cells.Item(157, 7).Value  = 100104
$ws.Cells.Item(157, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(157, 9).Value  = 100104003
$ws.Cells.Item(157, 10).Value = "Membrillo"
$ws.Cells.Item(157, 11).Value = "Champion"
$ws.Cells.Item(157, 12).Value = "Primera"
$ws.Cells.Item(157, 13).Value = 65
$ws.Cells.Item(157, 14).Value = 10000
$ws.Cells.Item(157, 15).Value = 10000
$ws.Cells.Item(157, 16).Value = 10000
$ws.Cells.Item(157, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(157, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(157, 19).Value = 556
$ws.Cells.Item(157, 20).Value = 18
